# Apply REST design updates: rename /admin/users endpoints to /users,
# add PUT support for the user endpoint, and document the access rule
# for who may call it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename endpoints (row 18 = collection endpoint, row 19 = item endpoint)
$ws.Range("D18").Value = "/users"
$ws.Range("D19").Value = "/users/{id}"

# /users/{id} now also supports PUT (update existing user), not just GET
$ws.Range("J19").Value = "GET/PUT"

# Document the access-control rule in column K, merged across the three
# /users rows (18-20), centered both horizontally and vertically.
$ws.Range("K18:K20").Merge()
$ws.Range("K18").Value = "only admin can access all users but user can only access his account"
$ws.Range("K18:K20").HorizontalAlignment = -4108
$ws.Range("K18:K20").VerticalAlignment = -4108

# Widen column K so the new note text fits
$ws.Columns.Item(11).ColumnWidth = 62.5

# Leave the cursor where the edit happened
$ws.Range("J19").Select()
